$wb = $excel.ActiveWorkbook

$wsRec = $wb.Worksheets.Item("Recommandations")
$wsRecData = @(
  @('BRVM - SERVICES PUBLICS', 0, 4, 1047, 106.15, '🟡 Observer', '➖ Neutre'),
  @('BRVM - AUTRES SECTEURS', 0, 1, 669.2, 669.2, '🟡 Observer', '➖ Neutre'),
  @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 3, 506.63, 163.32, '🟡 Observer', '➖ Neutre'),
  @('BRVM - DISTRIBUTION', 0, 1, 494.78, 494.78, '🟡 Observer', '➖ Neutre'),
  @('BRVM - SERVICES FINANCIERS', 0, 3, 436.47, 144.7, '🟡 Observer', '➖ Neutre'),
  @('BRVM - CONSOMMATION DE BASE    (**)', 0, 2, 431.6, 215.68, '🟡 Observer', '➖ Neutre'),
  @('BRVM-PRESTIGE', 0, 3, 429.14, 141.41, '🟡 Observer', '➖ Neutre'),
  @('BRVM - INDUSTRIELS', 0, 3, 395.5, 130.82, '🟡 Observer', '➖ Neutre'),
  @('BRVM - TRANSPORT', 0, 1, 370.83, 370.83, '🟡 Observer', '➖ Neutre'),
  @('BRVM - AGRICULTURE', 0, 1, 341.63, 341.63, '🟡 Observer', '➖ Neutre'),
  @('BRVM - ENERGIE', 0, 3, 330.85, 108.95, '🟡 Observer', '➖ Neutre'),
  @('BRVM - TELECOMMUNICATIONS', 0, 3, 281.01, 92.06, '🟡 Observer', '➖ Neutre'),
  @('BRVM – COMPOSITE TOTAL RETURN    (**)', 0, 2, 263.52, 130.87, '🟡 Observer', '➖ Neutre'),
  @('BRVM - INDUSTRIE  (**)', 0, 1, 257.22, 257.22, '🟡 Observer', '➖ Neutre'),
  @('BRVM-PRINCIPAL    (**)', 0, 1, 217.75, 217.75, '🟡 Observer', '➖ Neutre'),
  @('BRVM-PRINCIPAL  (**)', 0, 1, 217.65, 217.65, '🟡 Observer', '➖ Neutre'),
  @('BRVM-PRINCIPAL     (**)', 0, 1, 216.13, 216.13, '🟡 Observer', '➖ Neutre'),
  @('BRVM - CONSOMMATION DE BASE  (**)', 0, 1, 216.08, 216.08, '🟡 Observer', '➖ Neutre'),
  @('BRVM - FINANCES', 0, 1, 148.41, 148.41, '🟡 Observer', '➖ Neutre'),
  @('EVIOSYS PACKAGING SIEM CI (SEMC)', 2, 0, 14.47, 7.33, '🟡 Observer', '➖ Neutre'),
  @('ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)', 3, 0, 13.74, 2.08, '🟢 Achat', '✅ Renforcer'),
  @('SONATEL SN (SNTS)', 1, 0, 7.05, 7.05, '🟡 Observer', '➖ Neutre'),
  @('ORAGROUP TOGO (ORGT)', 2, 1, 5.77, 6.04, '🟡 Observer', '👀 À surveiller'),
  @('SAPH CI (SPHC)', 1, 0, 5.33, 5.33, '🟡 Observer', '➖ Neutre'),
  @('BANK OF AFRICA BF (BOABF)', 1, 0, 5.07, 5.07, '🟡 Observer', '➖ Neutre'),
  @('CIE CI (CIEC)', 1, 0, 4.66, 4.66, '🟡 Observer', '➖ Neutre'),
  @('NESTLE CI (NTLC)', 1, 0, 3.24, 3.24, '🟡 Observer', '➖ Neutre'),
  @('SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 1, 0, 3.23, 3.23, '🟡 Observer', '➖ Neutre'),
  @('ONATEL BF (ONTBF)', 1, 0, 1.22, 1.22, '🟡 Observer', '➖ Neutre'),
  @('SOLIBRA CI (SLBC)', 1, 0, 1.05, 1.05, '🟡 Observer', '➖ Neutre'),
  @('TOTAL', 0, 1, 0, 0, '🟡 Observer', '➖ Neutre'),
  @('BERNABE CI (BNBC)', 0, 1, -3.19, -3.19, '🟡 Observer', '➖ Neutre'),
  @('BANK OF AFRICA ML (BOAM)', 0, 1, -3.5, -3.5, '🟡 Observer', '➖ Neutre'),
  @('SAFCA CI (SAFC)', 0, 1, -3.93, -3.93, '🟡 Observer', '➖ Neutre'),
  @('ECOBANK TRANS. INCORP. TG (ETIT)', 0, 1, -4.35, -4.35, '🟡 Observer', '➖ Neutre'),
  @('FILTISAC CI (FTSC)', 0, 1, -5.24, -5.24, '🟡 Observer', '➖ Neutre'),
  @('SICOR CI (SICC)', 0, 1, -5.71, -5.71, '🟡 Observer', '➖ Neutre'),
  @('TOTALENERGIES MARKETING SN (TTLS)', 0, 1, -7.11, -7.11, '🟡 Observer', '➖ Neutre'),
  @('SETAO CI (STAC)', 0, 2, -10.09, -3.4, '🟡 Observer', '➖ Neutre'),
  @('CFAO MOTORS CI (CFAC)', 0, 3, -10.28, -6.45, '🔴 Vente', '⚠️ Risque de décrochage'),
  @('NEI-CEDA CI (NEIC)', 0, 2, -14.71, -7.21, '🟡 Observer', '➖ Neutre')
)
for ($i = 0; $i -lt $wsRecData.Count; $i++) {
    $row = $wsRecData[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $wsRec.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

$wsTop = $wb.Worksheets.Item("Top_YTD")
$wsTopData = @(
  @('BRVM - CONSOMMATION DE BASE    (**)', 897.3),
  @('BRVM - CONSOMMATION DISCRETIONNAIRE', 611.73),
  @('BRVM - SERVICES FINANCIERS', 501.74),
  @('BRVM-PRESTIGE', 487.79),
  @('BRVM – COMPOSITE TOTAL RETURN    (**)', 437.12),
  @('BRVM - INDUSTRIELS', 436.31),
  @('BRVM - ENERGIE', 340.86),
  @('BRVM - SERVICES PUBLICS', 327.12),
  @('BRVM - TELECOMMUNICATIONS', 272.69),
  @('BRVM-PRINCIPAL    (**)', 217.75)
)
for ($i = 0; $i -lt $wsTopData.Count; $i++) {
    $row = $wsTopData[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $wsTop.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
